$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename the "Chain (1ft)" item (row 6, col A) to " 04B Chain (1ft)"
#    and strip the explicit center-alignment style from its Link cell
#    (row 6, col B) so it falls back to the default/general style.
# ------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value2 = " 04B Chain (1ft)"
$ws.Cells.Item(6, 2).Style = "Normal"

# ------------------------------------------------------------------
# 2) Add a new line item for the Bowden Tube in row 20.
# ------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value2 = "Bowden Tube"
$ws.Cells.Item(20, 2).Value2 = "https://www.amazon.com/Jagwire-Sport-Housing-Slick-Lube-Titanium/dp/B085NBZMJS/ref=sr_1_1?crid=1A5WX5ADQYM0Y&dib=eyJ2IjoiMSJ9.CWV7EelBoN67bHKqG_VMGaAHwKq3lWAqkaCYopLdT43GyH4CDNeyWoQ_bFV_YrbZmhmwmsofP69GRzCCWYW_ULIkinZgZrdky8EGo_FPRa2GDLIPcrjwSu8T1nDFsZ03wyuLLatilRsdmpFkqvecV8S7AOhy1XjGVW6Ztcl1kgVL4_2zjOWbooP9z_kj4elJxdrMd7yL-uhr5ZcYm6F_Z725qnNy_c0-wNfaDt_xxhvLSabJzp2Ta9HzGpozVxtWGMo89NRT0qD_9iY1NxHHBedJp86w9TXZ0OfMm9s45WQ.7sQYhAOJzDOGkHvrSoGgv6F0CIzum6_7bjb1Zjv8zzc&dib_tag=se&keywords=Jagwire+Brake+Housing+CGX-SL+Slick-Lube+5+mm+%2810+m%29&qid=1726249759&s=sporting-goods&sprefix=jagwire+brake+housing+cgx-sl+slick-lube+5+mm+10+m+%2Csporting%2C99&sr=1-1"
$ws.Cells.Item(20, 3).Value2 = 1
$ws.Cells.Item(20, 4).Value2 = 61.97

# Center-align the new row's Item/Count/Price cells like the rest of the table
$ws.Cells.Item(20, 1).HorizontalAlignment = -4108
$ws.Cells.Item(20, 3).HorizontalAlignment = -4108
$ws.Cells.Item(20, 4).HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 3) Update the active selection to the newly added row, matching
#    where the author's cursor ended up after the edit.
# ------------------------------------------------------------------
[void]$ws.Range("A20:D20").Select()
